$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.457.99'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.861.02'
$ws.Range("E3").Value = '  +1.01%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.56'
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4776'
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3807'
$ws.Range("E8").Value = '  +3.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07313'
$ws.Range("E9").Value = '  +1.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9318'
$ws.Range("E10").Value = '  +0.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.74'
$ws.Range("E11").Value = '  +5.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07793'
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.864.93'
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.448'
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.550'
$ws.Range("E15").Value = '  +1.75%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.26'
$ws.Range("E16").Value = '  +1.64%  '
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008823'
$ws.Range("E18").Value = '  +1.94%  '
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.515.12'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.64'
$ws.Range("E21").Value = '  +0.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.097'
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.68'
$ws.Range("E23").Value = '  +0.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.944'
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.97'
$ws.Range("E25").Value = '  +1.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.46'
$ws.Range("E26").Value = '  +1.58%  '
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '115.53'
$ws.Range("E28").Value = '  +1.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.948'
$ws.Range("E29").Value = '  -0.32%  '
$ws.Range("E30").Value = '  +0.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.332'
$ws.Range("E31").Value = '  +1.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.204'
$ws.Range("E32").Value = '  +2.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7542'
$ws.Range("E33").Value = '  +2.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.585'
$ws.Range("E34").Value = '  +2.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.701'
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  +1.36%  '
$ws.Range("E37").Value = '  +4.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5558'
$ws.Range("E38").Value = '  +6.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05274'
$ws.Range("E39").Value = '  +0.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.989'
$ws.Range("E40").Value = '  +0.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.030'
$ws.Range("E41").Value = '  +0.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.595'
$ws.Range("E42").Value = '  +3.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1521'
$ws.Range("E43").Value = '  +0.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4870'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.67'
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.664'
$ws.Range("E47").Value = '  +3.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.93'
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '67.47'
$ws.Range("E49").Value = '  +2.60%  '
$ws.Range("E50").Value = '  +0.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9142'
$ws.Range("E51").Value = '  +3.12%  '
